$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) / E (Volume 1h) updates ---
$ws.Range("D2").Value = "37.501.02"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.076.30"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D5").Value = "235.38"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "57.51"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "0.393"
$ws.Range("D10").Value = "0.0779"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "14.39"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "20.75"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "0.783"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "5.21"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "2.069.29"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "37.395.99"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.70"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "226.73"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("D26").Value = "168.63"
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("D29").Value = "1.42"
$ws.Range("E29").Value = "  -6.11%  "
$ws.Range("D30").Value = "19.15"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "4.58"
$ws.Range("E32").Value = "  +1.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "5.62"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D43").Value = "97.19"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "4.22"
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").Value = "15.28"
$ws.Range("E48").Value = "  -3.24%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "2.263.72"
$ws.Range("E51").Value = "  +0.36%  "

# --- Row 27/28 swap: Kaspa <-> Cosmos (re-ranked) ---
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.134"
$ws.Range("E28").Value = "  +5.09%  "

# --- Row 41/42 swap: Maker <-> Cronos (re-ranked) ---
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "0.0957"
$ws.Range("E41").Value = "  -0.61%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.490.61"
$ws.Range("E42").Value = "  +2.55%  "

